$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "irr" table (left block, columns C/D/E) ---------------------------
# Completed items move out of the active list; new/ongoing items shift up.
$ws.Range("C4").Value = "get front end from monitoring and load it to git"
$ws.Range("C5").Value = "rebuild front end with most recent packages / add to repo"
$ws.Range("C6").Value = "redo irr for mortgages"
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()

# --- "rxr" table (right block, column G) --------------------------------
$ws.Range("G4").Value = "fix ratesheet table (remove useless space)"
$ws.Range("G5").Value = "fix sorting on ratesheet table"
$ws.Range("G6").Value = "create services/product segment"
$ws.Range("G7").Value = "redesign front end using bootstrap potentially"
$ws.Range("G8").Value = "download latest repo and have them run locaclly to make changes"
$ws.Range("G9").Value = "figure out how to load repo's to digital ocean and consider CI/CD"
$ws.Range("G10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("G12").ClearContents()
$ws.Range("G13").ClearContents()

# --- New "completed" section (rows 18-28) --------------------------------
# Bold section headers above each formerly-active list.
$ws.Range("C19").Value = "get upor/ucor for 360 mortgage loans"

$ws.Range("C18").Value = "completed"
$ws.Range("C18").Font.Bold = $true

$ws.Range("G19").Value = "login/registration fix"
$ws.Range("G19").Font.Color = 0
$ws.Range("G20").Value = "move dash board that is on front end to back end (as a copy)"
$ws.Range("G20").Font.Color = 0
$ws.Range("G21").Value = "update filters file for minimal operationability"
$ws.Range("G21").Font.Color = 0
$ws.Range("G22").Value = "fix upfront costs"
$ws.Range("G22").Font.Color = 0
$ws.Range("G23").Value = "fix apr calculation"
$ws.Range("G23").Font.Color = 0
$ws.Range("G24").Value = "create a pricing page"
$ws.Range("G24").Font.Color = 0
$ws.Range("G25").Value = "figure out api logic (renaming / combining / feeder)"
$ws.Range("G25").Font.Color = 0
$ws.Range("G26").Value = "redo design on the user 'logged in' page"
$ws.Range("G26").Font.Color = 0
$ws.Range("G27").Value = "my filter on button needs an alert that it may take a few seconds to load"
$ws.Range("G27").Font.Color = 0
$ws.Range("G28").Value = "create an api module that would generate the url"
$ws.Range("G28").Font.Color = 0

$ws.Range("G18").Value = "completed"
$ws.Range("G18").Font.Bold = $true

# --- Column width tweak (C widened, no longer auto best-fit) ------------
$ws.Columns.Item(3).ColumnWidth = 53.1640625

# --- Selection / active cell, matching the saved view --------------------
$ws.Range("C28").Select()
